$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.201.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3
$ws.Range("D3").Value = "'1.826.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("D4").Value = "'0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'241.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").Value = "'0.6200"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.41%  "

# Row 7
$ws.Range("D7").Value = "'0.9990"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.07355"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.39%  "

# Row 9
$ws.Range("D9").Value = "'0.2902"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "

# Row 10
$ws.Range("D10").Value = "'23.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "

# Row 11
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "'1.819.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13
$ws.Range("D13").Value = "'4.959"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.29%  "

# Row 14
$ws.Range("D14").Value = "'0.6637"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "

# Row 15
$ws.Range("D15").Value = "'82.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.92%  "

# Row 16
$ws.Range("D16").Value = "'0.000008940"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.79%  "

# Row 17
$ws.Range("D17").Value = "'5.850"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.48%  "

# Row 18
$ws.Range("D18").Value = "'29.164.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "

# Row 19
$ws.Range("D19").Value = "'2.067.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").Value = "'238.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.57%  "

# Row 21
$ws.Range("D21").Value = "'12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "

# Row 22
$ws.Range("D22").Value = "'0.9990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").Value = "'7.193"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.64%  "

# Row 24
$ws.Range("D24").Value = "'0.9995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").Value = "'158.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

# Row 26
$ws.Range("D26").Value = "'0.1418"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.42%  "

# Row 27
$ws.Range("D27").Value = "'8.438"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "

# Row 28
$ws.Range("D28").Value = "'17.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.09%  "

# Row 29
$ws.Range("D29").Value = "'1.486"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30
$ws.Range("D30").Value = "'0.05570"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.30%  "

# Row 31
$ws.Range("D31").Value = "'4.097"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "

# Row 32
$ws.Range("D32").Value = "'4.106"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.22%  "

# Row 33
$ws.Range("D33").Value = "'1.210"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "

# Row 34
$ws.Range("D34").Value = "'1.827"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "

# Row 35
$ws.Range("D35").Value = "'0.7364"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "

# Row 36
$ws.Range("E36").Value = "  -0.59%  "

# Row 37
$ws.Range("D37").Value = "'2.612"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.99%  "

# Row 38
$ws.Range("D38").Value = "'2.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.85%  "

# Row 39
$ws.Range("D39").Value = "'1.215.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.22%  "

# Row 40
$ws.Range("D40").Value = "'0.01759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.15%  "

# Row 41
$ws.Range("D41").Value = "'6.334"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.45%  "

# Row 42
$ws.Range("D42").Value = "'0.9203"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.23%  "

# Row 43
$ws.Range("D43").Value = "'0.9990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'101.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.39%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "'1.971.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'64.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.38%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.94%  "

# Row 48
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").Value = "'0.4018"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "

# Row 50
$ws.Range("D50").Value = "'9.096"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
$ws.Range("D51").Value = "'0.05769"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
